# Updates the "想去人数" (want-to-go count) values in column F across the
# four worksheets of the 北京-漫展信息 workbook, matching the regenerated
# data snapshot published to gh-pages (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value  = 359
$ws.Range("F7").Value  = 1173
$ws.Range("F9").Value  = 7111
$ws.Range("F11").Value = 88
$ws.Range("F12").Value = 2044
$ws.Range("F13").Value = 7975
$ws.Range("F15").Value = 52
$ws.Range("F16").Value = 5507
$ws.Range("F18").Value = 2410
$ws.Range("F19").Value = 1024
$ws.Range("F21").Value = 299
$ws.Range("F25").Value = 378
$ws.Range("F26").Value = 256
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 2347
$ws.Range("F34").Value = 6
$ws.Range("F36").Value = 1491
$ws.Range("F37").Value = 32
$ws.Range("F39").Value = 2314
$ws.Range("F42").Value = 7

# ---- 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 93
$ws.Range("F4").Value = 64
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 26

# ---- 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 254

# ---- 全部类型 (All types, combined view) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 254
$ws.Range("F7").Value  = 93
$ws.Range("F8").Value  = 359
$ws.Range("F9").Value  = 1173
$ws.Range("F11").Value = 7111
$ws.Range("F13").Value = 88
$ws.Range("F14").Value = 2044
$ws.Range("F15").Value = 7975
$ws.Range("F17").Value = 52
$ws.Range("F18").Value = 5507
$ws.Range("F20").Value = 2410
$ws.Range("F21").Value = 1024
$ws.Range("F27").Value = 64
$ws.Range("F28").Value = 378
$ws.Range("F29").Value = 9
$ws.Range("F30").Value = 2347
$ws.Range("F35").Value = 5
$ws.Range("F37").Value = 6
$ws.Range("F39").Value = 26
$ws.Range("F40").Value = 1491
$ws.Range("F41").Value = 32
$ws.Range("F43").Value = 2314
$ws.Range("F47").Value = 7
